# Word COM-interop script implementing the diff:
# "The Inseparable Nexus of Science and Human Endeavor" (Benjamin Rossner)
#  -> "The Enchanting World of Biology: Unveiling the Secrets of Life" (Dr. Henrietta Adams)

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $f = $d.Content.Find
    $f.ClearFormatting()
    $f.Text = ""
    $ok = $f.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false)
    if ($ok) {
        $f.Parent.Text = $replace
    }
    return $ok
}

# ---- Title ----
Replace-Text "The Inseparable Nexus of Science and Human Endeavor" "The Enchanting World of Biology: Unveiling the Secrets of Life"

# ---- Author name: "Benjamin Rossner" -> "Dr" + "." + " Henrietta Adams" (3 runs) ----
$f = $d.Content.Find
$f.ClearFormatting()
$f.Text = ""
$f.Execute("Benjamin Rossner", $true, $false, $false, $false, $false, $true, 1, $false)
$f.Parent.Text = "Dr"
$endPos = $f.Parent.End
$r1 = $d.Range($endPos, $endPos)
$r1.InsertAfter(".")
$endPos2 = $endPos + 1
$r2 = $d.Range($endPos2, $endPos2)
$r2.InsertAfter(" Henrietta Adams")

# ---- Email ----
Replace-Text "benjamin" "henrietta"
Replace-Text "rossner@gmail" "adams@kingsley"
Replace-Text "com" "edu"

# ---- Body paragraph 1st block ----
Replace-Text "Science represents not simply a body of knowledge but an intellectual voyage, transcending the boundaries of disciplines and entwining with the fabric of human experience" "Biology, the study of life, is a fascinating field that invites us to unravel the enigmatic secrets of existence"

Replace-Text " Through meticulous observation and exploration, science unravels the intricate workings of the cosmos, penetrating the depths of understanding to unearth fundamental truths" " From the intricate workings of cells to the majestic diversity of organisms that grace our planet, biology captivates our imagination and inspires awe"

Replace-Text " It unveils the secrets of our physical world, providing insights into the composition of matter, the motion of celestial bodies, and the forces that shape the universe" " It probes the fundamental questions about the origin and nature of life, exploring the exquisite intricacies of living systems and their interactions within the natural world"

Replace-Text " Science enhances our comprehension of life itself, elucidating the mysteries of biology, chemistry, and genetics, revealing the intricate interplay of cells and organisms" " Biology not only enriches our understanding of the world but also holds the potential to address urgent challenges, from disease prevention and environmental conservation to the quest for sustainable living"

# New runs inserted after that run (before the existing "." run)
$f = $d.Content.Find
$f.ClearFormatting()
$f.Text = ""
$f.Execute("quest for sustainable living", $true, $false, $false, $false, $false, $true, 1, $false)
$endPos = $f.Parent.End
$r1 = $d.Range($endPos, $endPos)
$r1.InsertAfter(".")
$endPos2 = $endPos + 1
$r2 = $d.Range($endPos2, $endPos2)
$r2.InsertAfter(" Join us on an enthralling journey as we delve into the wonders of biology, unlocking the mysteries of life's grand tapestry")

Replace-Text "Yet, science extends beyond its methodical exploration of the natural world" "Biology's grand tapestry unveils the remarkable diversity of life on Earth, showcasing the intricate adaptations and interconnections among organisms"

Replace-Text " Its tendrils reach into the realm of human endeavor, intertwining with our social, cultural, and philosophical pursuits" " From the microscopic realm teeming with microbes to the towering giants of the plant kingdom, each species occupies a unique niche, playing a vital role in maintaining the delicate balance of ecosystems"

Replace-Text " Science sculpts our perception of reality, challenges age-old beliefs, and influences the way we interact with our surroundings" " Evolution, the driving force behind this diversity, has shaped the remarkable complexity and resilience of life, crafting organisms capable of thriving in a myriad of environments, from scorching deserts to frigid polar regions"

# Collapse 3 runs (" It empowers...", ".", " Moreover, science...") into one replaced run + delete the other two
Replace-Text " It empowers us to mitigate diseases, harness energy, and traverse vast distances" " Discover the awe-inspiring beauty of nature's designs and the intricate interrelationships that sustain the web of life"

# delete the "." run immediately followed by " Moreover, science..." run (together, exact text span)
$f = $d.Content.Find
$f.ClearFormatting()
$f.Text = ""
$f.Execute(". Moreover, science serves as a catalyst for innovation, spurring technological advancements that transform industries and redefine the limits of human potential", $true, $false, $false, $false, $false, $true, 1, $false)
$f.Parent.Text = ""

Replace-Text "In recent times, the symbiotic relationship between science and human endeavor has become increasingly apparent" "Biology is not merely a theoretical pursuit; it has tangible implications for our daily lives"

Replace-Text " Globalization has accelerated the sharing of scientific knowledge and expertise, fostering unprecedented collaboration among researchers worldwide" " Advances in medical research have produced life-saving vaccines and therapies, significantly improving human health"

Replace-Text " Interdisciplinary approaches to problem-solving have yielded transformative insights and breakthroughs" " Innovations in agriculture have boosted crop yields, ensuring food security for a growing population"

Replace-Text " Science has become a cornerstone of public policy, shaping decisions that impact healthcare, environmental conservation, and economic development" " Our understanding of genetics influences advances in DNA analysis, enabling advancements in forensics and personalized medicine"

# New runs: "." + " Biology also plays a crucial role..." + "." + " Engage with biology and " + (lastRenderedPageBreak) "contribute to these endeavors..."
$f = $d.Content.Find
$f.ClearFormatting()
$f.Text = ""
$f.Execute("enabling advancements in forensics and personalized medicine", $true, $false, $false, $false, $false, $true, 1, $false)
$endPos = $f.Parent.End
$r1 = $d.Range($endPos, $endPos)
$r1.InsertAfter(".")
$p1 = $endPos + 1
$r2 = $d.Range($p1, $p1)
$r2.InsertAfter(" Biology also plays a crucial role in addressing global challenges like climate change and pollution, as scientists strive to create sustainable solutions for a healthier planet")
$p2 = $r2.End
$r3 = $d.Range($p2, $p2)
$r3.InsertAfter(".")
$p3 = $r3.End
$r4 = $d.Range($p3, $p3)
$r4.InsertAfter(" Engage with biology and ")
$p4 = $r4.End
$r5 = $d.Range($p4, $p4)
$r5.InsertAfter("contribute to these endeavors, leaving a positive impact on the world through scientific discovery")
